# Add "Day 6" entry as a new row (row 7) at the bottom of the tracker table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting from the row above (B6) so the new date cell
# (B7) reuses the existing date cell style instead of minting a new one.
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A7").Value = "Day 6"
$ws.Range("B7").Value = 45808         # 5/31/2025
$ws.Range("D7").Value = "Remove Element"
$ws.Range("C7").Value = "Longest Consecutive Sequence"
$ws.Range("E7").Value = "Search a 2D Matrix"
$ws.Range("F7").Value = "In-place Array, Hashing, 2D Binary Search"
$ws.Range("G7").Value = "S"
$ws.Range("H7").Value = "YES"

$ws.Range("H8").Select()
